$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C10 holds the "Integer min" lower bound for rule R30 (row 10).
# Restore it from 18 back to 1, per the target revision.
$ws.Range("C10").Value = 1
